$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.847.11'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '2.400.75'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').Value = '  +0.51%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '565.25'
$c.ClearFormats()
$ws.Range('E5').Value = '  -0.03%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '141.90'
$c.ClearFormats()
$ws.Range('E6').Value = '  +3.00%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range('E7').Value = '  -0.45%  '
$ws.Range('E8').Value = '  +2.76%  '
$ws.Range('D9').Value = '2.407.87'
$ws.Range('E9').Value = '  +1.06%  '
$ws.Range('E10').Value = '  +2.14%  '
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('E12').Value = '  +3.29%  '
$ws.Range('E13').Value = '  +3.46%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '26.40'
$c.ClearFormats()
$ws.Range('E14').Value = '  +2.67%  '
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').Value = '2.837.08'
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('D17').Value = '60.680.95'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '2.404.55'
$ws.Range('E18').Value = '  +0.23%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '8.04'
$c.ClearFormats()
$ws.Range('E19').Value = '  +3.90%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '10.70'
$c.ClearFormats()
$ws.Range('E20').Value = '  +1.70%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '324.20'
$c.ClearFormats()
$ws.Range('E21').Value = '  +1.15%  '
$ws.Range('E22').Value = '  +2.27%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '6.04'
$c.ClearFormats()
$ws.Range('E23').Value = '  -0.51%  '
$ws.Range('E24').Value = '  -0.24%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '1.89'
$c.ClearFormats()
$ws.Range('E25').Value = '  +5.38%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '65.04'
$c.ClearFormats()
$ws.Range('E26').Value = '  +0.54%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '586.95'
$c.ClearFormats()
$ws.Range('E27').Value = '  +2.77%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '8.20'
$c.ClearFormats()
$ws.Range('E28').Value = '  +0.55%  '
$ws.Range('D29').Value = '0.0₃0942'
$ws.Range('E29').Value = '  +3.71%  '
$ws.Range('E30').Value = '  -1.07%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '8.02'
$c.ClearFormats()
$ws.Range('E31').Value = '  +3.39%  '
$ws.Range('E32').Value = '  +1.78%  '
$ws.Range('E33').Value = '  +0.63%  '
$ws.Range('E34').Value = '  +1.56%  '
$ws.Range('E35').Value = '  +6.68%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.ClearFormats()
$ws.Range('E36').Value = '  -0.61%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '152.48'
$c.ClearFormats()
$ws.Range('E37').Value = '  +0.71%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.371'
$c.ClearFormats()
$ws.Range('E38').Value = '  +2.25%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '4.62'
$c.ClearFormats()
$ws.Range('E39').Value = '  +1.84%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '18.32'
$c.ClearFormats()
$ws.Range('E40').Value = '  +1.53%  '
$ws.Range('E41').Value = '  +3.19%  '
$ws.Range('E42').Value = '  -0.12%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.52'
$c.ClearFormats()
$ws.Range('E43').Value = '  +12.76%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '1.68'
$c.ClearFormats()
$ws.Range('E44').Value = '  +2.42%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '41.60'
$c.ClearFormats()
$ws.Range('E45').Value = '  +1.29%  '
$ws.Range('D46').Value = '0.0₆0278'
$ws.Range('E46').Value = '  +7.30%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '141.85'
$c.ClearFormats()
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('E48').Value = '  +1.60%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.591'
$c.ClearFormats()
$ws.Range('E49').Value = '  +1.81%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0509'
$c.ClearFormats()
$ws.Range('E50').Value = '  +2.70%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '19.48'
$c.ClearFormats()
$ws.Range('E51').Value = '  +2.16%  '
